$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. New risk entry "RISK 4" recorded in the Impact/Probability matrix (cell D3),
#    which becomes a new shared string (sharedStrings.xml gains "RISK 4").
$ws.Range("D3").Value = "RISK 4"

# 2. The "Risks:" text box (TextBox 4 / shape id 5) gets a 4th bullet describing
#    the new risk: "4. Training". The COM text-box API only exposes whole-text
#    replacement (Characters().Text), so rebuild the full label with the extra
#    line appended after item 3, keeping the trailing blank line that the
#    original text box ended with.
$risksBox = $ws.Shapes.Item("TextBox 4")
$lines = @(
    "Risks:",
    "1. No Order Data",
    "2. No Factory Data",
    "3. Factory Locator Subsystem Server Down (centralized architecture)",
    "4. Training",
    ""
)
$risksBox.TextFrame.Characters().Text = [string]::Join([char]13, $lines)

# 3. Cosmetic: the saved workbook's selection moved to K3.
[void]$ws.Range("K3").Select()
